$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 528; this shifts every subsequent row down by
# one (old row 528 -> 529, ..., old row 606 -> 607) and extends the used
# range to A1:R607.
$ws.Rows.Item(528).Insert()

# Populate the newly inserted row 528 with the new price-report record.
$ws.Range("A528").Value = 9
$ws.Range("B528").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C528").Value = "Metropolitana"
$ws.Range("D528").Value = 45131
$ws.Range("E528").Value = 13
$ws.Range("F528").Value = 100112052
$ws.Range("G528").Value = "Albahaca"
$ws.Range("H528").Value = "Sin especificar"
$ws.Range("I528").Value = "Primera"
$ws.Range("J528").Value = 430
$ws.Range("K528").Value = 4000
$ws.Range("L528").Value = 4500
$ws.Range("M528").Value = 4250
$ws.Range("N528").Value = "$/paquete"
$ws.Range("O528").Value = "Región de Arica y Parinacota"
$ws.Range("P528").Value = 4250
$ws.Range("Q528").Value = 1
$ws.Range("R528").Value = "Hortaliza"
